$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B ("Right Index") ---
# The static values in B2:B42 (which were simply "A + 100") are replaced
# by a live formula "=A + 41". Assign the first formula directly and then
# fill the rest of the range so Excel records it as one shared formula
# group, exactly as in the authored workbook (f t="shared" si="0").
$ws.Range("B2").Formula = "=A2+41"
$ws.Range("B3:B42").Formula = "=A3+41"

# --- B42 formatting ---
# Row 42 previously used a distinct "final row" cell style for column B
# (s=11). After the edit it uses the same style as the rest of the
# column (s=1, the style already applied to B2:B41), so copy that
# formatting (not the value/formula) from B41 onto B42.
$ws.Range("B41").Copy() | Out-Null
$ws.Range("B42").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# The last row grows slightly taller (24 -> 24.75) once its formatting
# lines up with the rest of the table.
$ws.Rows.Item(42).RowHeight = 24.75

# --- Selection ---
# Move the active selection to where the author finished editing.
$ws.Range("L40").Select() | Out-Null
